$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 12.2
$ws.Range("I11").Value = 12.2
$ws.Range("K11").Value = 12.2
$ws.Range("M11").Value = 127.8
$ws.Range("H76").Value = 4034.75
$ws.Range("I76").Value = 3897.5
$ws.Range("K76").Value = 3897.5
$ws.Range("M76").Value = -3582.5
$ws.Range("H79").Value = 4034.75
$ws.Range("I79").Value = 3897.5
$ws.Range("K79").Value = 3897.5
$ws.Range("M79").Value = -2805.5
$ws.Range("H97").Value = 716047.2
$ws.Range("J97").Value = 716047.2
$ws.Range("L97").Value = 2148141.6
$ws.Range("N97").Value = -2149133.6
$ws.Range("H127").Value = 62502800
$ws.Range("I127").Value = 76924296
$ws.Range("K127").Value = 230772888
$ws.Range("M127").Value = -230767928
$ws.Range("H131").Value = 2969.2222
$ws.Range("I131").Value = 3203.625
$ws.Range("K131").Value = 9610.875
$ws.Range("M131").Value = -4570.875
$ws.Range("H135").Value = 952.2
$ws.Range("I135").Value = 952.2
$ws.Range("K135").Value = 8569.800000000001
$ws.Range("M135").Value = -6034.800000000001
$ws.Range("H137").Value = 2710.963
$ws.Range("I137").Value = 2724.4412
$ws.Range("J137").Value = 2688.05
$ws.Range("K137").Value = 8173.323600000001
$ws.Range("L137").Value = 8064.150000000001
$ws.Range("M137").Value = -5623.323600000001
$ws.Range("N137").Value = -13164.15
$ws.Range("H141").Value = 1583.7354
$ws.Range("I141").Value = 512.8889
$ws.Range("K141").Value = 1538.6667
$ws.Range("M141").Value = 3641.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8803.548000000001
$ws.Range("I32").Value = 6975.364
$ws.Range("J32").Value = 26040.715
$ws.Range("K32").Value = 6975.364
$ws.Range("L32").Value = 26040.715
$ws.Range("M32").Value = -6688.364
$ws.Range("N32").Value = -26614.715
$ws.Range("H74").Value = 2873.628
$ws.Range("I74").Value = 2436.8157
$ws.Range("J74").Value = 6193.4
$ws.Range("K74").Value = 2436.8157
$ws.Range("L74").Value = 6193.4
$ws.Range("M74").Value = -1562.8157
$ws.Range("N74").Value = -7941.4
$ws.Range("H77").Value = 2873.628
$ws.Range("I77").Value = 2436.8157
$ws.Range("J77").Value = 6193.4
$ws.Range("K77").Value = 12184.0785
$ws.Range("L77").Value = 30967
$ws.Range("M77").Value = -7816.0785
$ws.Range("N77").Value = -39703
$ws.Range("H97").Value = 1299.68
$ws.Range("I97").Value = 935.05
$ws.Range("J97").Value = 2758.2
$ws.Range("K97").Value = 935.05
$ws.Range("L97").Value = 2758.2
$ws.Range("M97").Value = -439.05
$ws.Range("N97").Value = -3750.2
$ws.Range("H122").Value = 6354.909
$ws.Range("I122").Value = 2705.4443
$ws.Range("K122").Value = 8116.3329
$ws.Range("M122").Value = -5666.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1932.6666
$ws.Range("I94").Value = 2699
$ws.Range("K94").Value = 2699
$ws.Range("M94").Value = -2248
$ws.Range("H134").Value = 2444.0557
$ws.Range("I134").Value = 1912
$ws.Range("J134").Value = 4040.2222
$ws.Range("K134").Value = 5736
$ws.Range("L134").Value = 12120.6666
$ws.Range("M134").Value = -3201
$ws.Range("N134").Value = -17190.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2144.3713
$ws.Range("I31").Value = 2080.3333
$ws.Range("J31").Value = 2360.5
$ws.Range("K31").Value = 2080.3333
$ws.Range("L31").Value = 2360.5
$ws.Range("M31").Value = -1785.3333
$ws.Range("N31").Value = -2950.5
$ws.Range("H34").Value = 2144.3713
$ws.Range("I34").Value = 2080.3333
$ws.Range("J34").Value = 2360.5
$ws.Range("K34").Value = 2080.3333
$ws.Range("L34").Value = 2360.5
$ws.Range("M34").Value = -1878.3333
$ws.Range("N34").Value = -2764.5
$ws.Range("H134").Value = 1002.0548
$ws.Range("I134").Value = 950.18335
$ws.Range("K134").Value = 2850.55005
$ws.Range("M134").Value = -315.5500499999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1016
$ws.Range("J45").Value = 1016
$ws.Range("L45").Value = 3048
$ws.Range("N45").Value = -4112
$ws.Range("H68").Value = 5210481
$ws.Range("J68").Value = 4650.3335
$ws.Range("L68").Value = 13951.0005
$ws.Range("N68").Value = -15573.0005
$ws.Range("H71").Value = 5210481
$ws.Range("J71").Value = 4650.3335
$ws.Range("L71").Value = 41853.0015
$ws.Range("N71").Value = -49965.0015
$ws.Range("H122").Value = 681
$ws.Range("J122").Value = 703
$ws.Range("L122").Value = 6327
$ws.Range("N122").Value = -11227

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 11633.333
$ws.Range("J92").Value = 11633.333
$ws.Range("L92").Value = 11633.333
$ws.Range("N92").Value = -15377.333
$ws.Range("H97").Value = 2293.5557
$ws.Range("I97").Value = 2007.5518
$ws.Range("K97").Value = 2007.5518
$ws.Range("M97").Value = -1511.5518
$ws.Range("H126").Value = 6805
$ws.Range("I126").Value = 8866.666999999999
$ws.Range("J126").Value = 3712.5
$ws.Range("K126").Value = 26600.001
$ws.Range("L126").Value = 11137.5
$ws.Range("M126").Value = -24130.001
$ws.Range("N126").Value = -16077.5
$ws.Range("H132").Value = 7249.69
$ws.Range("I132").Value = 8000.373
$ws.Range("K132").Value = 24001.119
$ws.Range("M132").Value = -21471.119

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3533.4583
$ws.Range("I68").Value = 3617.318
$ws.Range("J68").Value = 2611
$ws.Range("K68").Value = 3617.318
$ws.Range("L68").Value = 2611
$ws.Range("M68").Value = -2868.318
# N68 is a new cell introduced on this row (previously absent)
$ws.Range("N68").Value = -4109
$ws.Range("H71").Value = 3533.4583
$ws.Range("I71").Value = 3617.318
$ws.Range("J71").Value = 2611
$ws.Range("K71").Value = 18086.59
$ws.Range("L71").Value = 13055
$ws.Range("M71").Value = -14342.59
# N71 is a new cell introduced on this row (previously absent)
$ws.Range("N71").Value = -20543
$ws.Range("H132").Value = 4146.3423
$ws.Range("I132").Value = 2776.3555
$ws.Range("K132").Value = 8329.066500000001
$ws.Range("M132").Value = -5799.066500000001
$ws.Range("H136").Value = 2501.625
$ws.Range("I136").Value = 1762.7273
$ws.Range("K136").Value = 5288.1819
$ws.Range("M136").Value = -2738.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1610.2394
$ws.Range("I132").Value = 1237.1754
$ws.Range("J132").Value = 3129.1428
$ws.Range("K132").Value = 3711.5262
$ws.Range("L132").Value = 9387.428400000001
$ws.Range("M132").Value = -1181.5262
$ws.Range("N132").Value = -14447.4284
